$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for existing data rows 2-433 from 45202 to 45203
$ws.Range("C2:C433").Value = 45203

# 2. Row 433 gains an explicit row height (15, customHeight)
$ws.Rows.Item(433).RowHeight = 15

# 3. Append a new data row (434) for case "A 46882-2023"
$ws.Range("A434").Value = "A 46882-2023"

$ws.Range("B434").Value = 45201
$ws.Range("B434").NumberFormat = $ws.Range("B433").NumberFormat()

$ws.Range("C434").Value = 45203
$ws.Range("C434").NumberFormat = $ws.Range("C433").NumberFormat()

$ws.Range("D434").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E434").Value = "STRÄNGNÄS"

$ws.Range("G434").Value = 8.5
$ws.Range("H434").Value = 0
$ws.Range("I434").Value = 0
$ws.Range("J434").Value = 0
$ws.Range("K434").Value = 0
$ws.Range("L434").Value = 0
$ws.Range("M434").Value = 0
$ws.Range("N434").Value = 0
$ws.Range("O434").Value = 0
$ws.Range("P434").Value = 0
$ws.Range("Q434").Value = 0

$ws.Range("R434").WrapText = $true
